# Update BP terminal gate pricing table with the latest effective-date prices.
# The former "current" day (2026-02-10, serial 46063) rows shift down into the
# "previous" day slot, new prices for 2026-02-11 (serial 46064) take the "current"
# day slot, and the oldest day (2026-02-07, serial 46060) is dropped from the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New South Wales
$ws.Cells.Item(8, 1).Value = 46064
$ws.Cells.Item(8, 4).Value = 159.09
$ws.Cells.Item(8, 5).Value = 148.81
$ws.Cells.Item(8, 6).Value = 158.81
$ws.Cells.Item(8, 7).Value = 148.69999999999999

$ws.Cells.Item(9, 1).Value = 46064
$ws.Cells.Item(9, 4).Value = 159.09
$ws.Cells.Item(9, 5).Value = 148.81
$ws.Cells.Item(9, 6).Value = 158.81
$ws.Cells.Item(9, 7).Value = 148.69999999999999

$ws.Cells.Item(10, 1).Value = 46064
$ws.Cells.Item(10, 4).Value = 160.62
$ws.Cells.Item(10, 5).Value = 151.32
$ws.Cells.Item(10, 6).Value = 161.32
$ws.Cells.Item(10, 7).Value = 151.56

$ws.Cells.Item(11, 1).Value = 46063
$ws.Cells.Item(11, 4).Value = 159.35
$ws.Cells.Item(11, 5).Value = 148.62
$ws.Cells.Item(11, 6).Value = 158.62
$ws.Cells.Item(11, 7).Value = 148.51

$ws.Cells.Item(12, 1).Value = 46063
$ws.Cells.Item(12, 4).Value = 159.35
$ws.Cells.Item(12, 5).Value = 148.62
$ws.Cells.Item(12, 6).Value = 158.62
$ws.Cells.Item(12, 7).Value = 148.51

$ws.Cells.Item(13, 1).Value = 46063
$ws.Cells.Item(13, 4).Value = 160.74
$ws.Cells.Item(13, 5).Value = 150.77000000000001
$ws.Cells.Item(13, 6).Value = 160.77000000000001
$ws.Cells.Item(13, 7).Value = 151.02000000000001

# Northern Territory
$ws.Cells.Item(17, 1).Value = 46064
$ws.Cells.Item(17, 4).Value = 164.85
$ws.Cells.Item(17, 5).Value = 154.74
$ws.Cells.Item(17, 6).Value = 164.74

$ws.Cells.Item(18, 1).Value = 46063
$ws.Cells.Item(18, 4).Value = 164.95
$ws.Cells.Item(18, 5).Value = 154.16999999999999
$ws.Cells.Item(18, 6).Value = 164.17

# Queensland
$ws.Cells.Item(22, 1).Value = 46064
$ws.Cells.Item(22, 4).Value = 160.27000000000001
$ws.Cells.Item(22, 5).Value = 150.93
$ws.Cells.Item(22, 6).Value = 160.53
$ws.Cells.Item(22, 7).Value = 152.68

$ws.Cells.Item(23, 1).Value = 46064
$ws.Cells.Item(23, 4).Value = 165.6
$ws.Cells.Item(23, 5).Value = 157.01
$ws.Cells.Item(23, 6).Value = 167.01

$ws.Cells.Item(24, 1).Value = 46064
$ws.Cells.Item(24, 4).Value = 165.77
$ws.Cells.Item(24, 5).Value = 157.6
$ws.Cells.Item(24, 6).Value = 167.6

$ws.Cells.Item(25, 1).Value = 46064
$ws.Cells.Item(25, 4).Value = 165.77
$ws.Cells.Item(25, 5).Value = 157.12
$ws.Cells.Item(25, 6).Value = 167.12
$ws.Cells.Item(25, 7).Value = 157.97999999999999

$ws.Cells.Item(26, 1).Value = 46064
$ws.Cells.Item(26, 4).Value = 165.38
$ws.Cells.Item(26, 5).Value = 158.72
$ws.Cells.Item(26, 6).Value = 168.72

$ws.Cells.Item(27, 1).Value = 46063
$ws.Cells.Item(27, 4).Value = 160.41999999999999
$ws.Cells.Item(27, 5).Value = 150.74
$ws.Cells.Item(27, 6).Value = 160.34
$ws.Cells.Item(27, 7).Value = 152.49

$ws.Cells.Item(28, 1).Value = 46063
$ws.Cells.Item(28, 4).Value = 165.73
$ws.Cells.Item(28, 5).Value = 156.44999999999999
$ws.Cells.Item(28, 6).Value = 166.45

$ws.Cells.Item(29, 1).Value = 46063
$ws.Cells.Item(29, 4).Value = 165.89
$ws.Cells.Item(29, 5).Value = 157.06
$ws.Cells.Item(29, 6).Value = 167.06

$ws.Cells.Item(30, 1).Value = 46063
$ws.Cells.Item(30, 4).Value = 165.89
$ws.Cells.Item(30, 5).Value = 156.59
$ws.Cells.Item(30, 6).Value = 166.59
$ws.Cells.Item(30, 7).Value = 157.44

$ws.Cells.Item(31, 1).Value = 46063
$ws.Cells.Item(31, 4).Value = 165.5
$ws.Cells.Item(31, 5).Value = 158.16999999999999
$ws.Cells.Item(31, 6).Value = 168.17

# South Australia
$ws.Cells.Item(35, 1).Value = 46064
$ws.Cells.Item(35, 4).Value = 159.11000000000001
$ws.Cells.Item(35, 5).Value = 149.01
$ws.Cells.Item(35, 6).Value = 158.01

$ws.Cells.Item(36, 1).Value = 46063
$ws.Cells.Item(36, 4).Value = 159.22999999999999
$ws.Cells.Item(36, 5).Value = 148.46
$ws.Cells.Item(36, 6).Value = 157.46

# Tasmania
$ws.Cells.Item(40, 1).Value = 46064
$ws.Cells.Item(40, 4).Value = 165.41
$ws.Cells.Item(40, 5).Value = 156.69
$ws.Cells.Item(40, 6).Value = 166.69

$ws.Cells.Item(41, 1).Value = 46064
$ws.Cells.Item(41, 4).Value = 165.13
$ws.Cells.Item(41, 5).Value = 157.12
$ws.Cells.Item(41, 6).Value = 167.12

$ws.Cells.Item(42, 1).Value = 46063
$ws.Cells.Item(42, 4).Value = 165.56
$ws.Cells.Item(42, 5).Value = 156.19
$ws.Cells.Item(42, 6).Value = 166.19

$ws.Cells.Item(43, 1).Value = 46063
$ws.Cells.Item(43, 4).Value = 165.27
$ws.Cells.Item(43, 5).Value = 156.61000000000001
$ws.Cells.Item(43, 6).Value = 166.61

# Victoria
$ws.Cells.Item(47, 1).Value = 46064
$ws.Cells.Item(47, 4).Value = 160
$ws.Cells.Item(47, 5).Value = 150.37
$ws.Cells.Item(47, 6).Value = 160.37

$ws.Cells.Item(48, 1).Value = 46064
$ws.Cells.Item(48, 4).Value = 159.66
$ws.Cells.Item(48, 5).Value = 150.32
$ws.Cells.Item(48, 6).Value = 160.32

$ws.Cells.Item(49, 1).Value = 46063
$ws.Cells.Item(49, 4).Value = 160.19
$ws.Cells.Item(49, 5).Value = 150.16999999999999
$ws.Cells.Item(49, 6).Value = 160.16999999999999

$ws.Cells.Item(50, 1).Value = 46063
$ws.Cells.Item(50, 4).Value = 159.85
$ws.Cells.Item(50, 5).Value = 150.13
$ws.Cells.Item(50, 6).Value = 160.13

# Western Australia
$ws.Cells.Item(54, 1).Value = 46064
$ws.Cells.Item(54, 4).Value = 174.51
$ws.Cells.Item(54, 5).Value = 164.4
$ws.Cells.Item(54, 6).Value = 174.4

$ws.Cells.Item(55, 1).Value = 46064
$ws.Cells.Item(55, 4).Value = 163.89
$ws.Cells.Item(55, 5).Value = 162.44999999999999
$ws.Cells.Item(55, 6).Value = 172.45

$ws.Cells.Item(56, 1).Value = 46064
$ws.Cells.Item(56, 4).Value = 163.76

$ws.Cells.Item(57, 1).Value = 46064
$ws.Cells.Item(57, 4).Value = 164.46
$ws.Cells.Item(57, 5).Value = 156.87

$ws.Cells.Item(58, 1).Value = 46064
$ws.Cells.Item(58, 4).Value = 160.22999999999999
$ws.Cells.Item(58, 5).Value = 152.77000000000001
$ws.Cells.Item(58, 6).Value = 162.77000000000001

$ws.Cells.Item(59, 1).Value = 46064
$ws.Cells.Item(59, 4).Value = 167.25
$ws.Cells.Item(59, 5).Value = 162.75

$ws.Cells.Item(60, 1).Value = 46063
$ws.Cells.Item(60, 4).Value = 174.64

$ws.Cells.Item(61, 1).Value = 46063
$ws.Cells.Item(61, 4).Value = 164.02
$ws.Cells.Item(61, 5).Value = 161.88
$ws.Cells.Item(61, 6).Value = 171.88

$ws.Cells.Item(62, 1).Value = 46063
$ws.Cells.Item(62, 4).Value = 163.89

$ws.Cells.Item(63, 1).Value = 46063
$ws.Cells.Item(63, 4).Value = 164.56
$ws.Cells.Item(63, 5).Value = 156.30000000000001

$ws.Cells.Item(64, 1).Value = 46063
$ws.Cells.Item(64, 4).Value = 160.33000000000001
$ws.Cells.Item(64, 5).Value = 152.19999999999999
$ws.Cells.Item(64, 6).Value = 162.19999999999999

$ws.Cells.Item(65, 1).Value = 46063
$ws.Cells.Item(65, 4).Value = 167.34
$ws.Cells.Item(65, 5).Value = 162.25
